# Update children's full names in column A of the active sheet.
# Mapping derived from old value -> new value, applied per-row so that
# only the intended cells are touched (some names share words, so we
# match on exact current value rather than blind row numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$renames = @{
    "Boden Williams"               = "Boden Nelson Williams"
    "Calum Carroll ward"           = "Calum Thomas Carroll Ward"
    "Charlotte Guyler"             = "Lottie Charlotte Guyler"
    "Charlotte Rose Doyle"         = "Charlotte Doyle"
    "Eleanor Wadden"               = "Eleanor Niamh Wadden"
    "Elsie Williams"                = "Elsie Nelson Williams"
    "Eve Otoole"                   = "Eve O Toole"
    "Jack Vickers McGerr"          = "JP Vickers McGerr"
    "Juno Luna Hynes Byrne"        = "Juno Hynes Byrne"
    "Katie Vickers Mc Gerr"        = "Katie Mcgerr"
    "Lily Grnik"                   = "Lily Gornik"
    "Lily Kathy May Corcoran"      = "Lily Corcoran"
    "Marc Aurele Gaaloul Donnelly" = "Marc Gaaloul Donnelly"
    "Naoise Siochr"                = "Naoise O Siochru"
    "Ray OCleirigh"                = "Ray O Cleirigh"
    "Theo OShaughnessy"            = "Theo O Shaughnessy"
    "Toms Hobbs"                   = "Toms Carmody Finnegan"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($null -ne $current -and $renames.ContainsKey([string]$current)) {
        $cell.Value = $renames[[string]$current]
    }
}
